$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (existing rows 2 and 3 shift down to 3 and 4)
$ws.Rows.Item(2).Insert()

# Match the formatting of the other data rows (Insert() pulls the header's
# bold/border format by default, so pull the plain data-row format down
# from what is now row 3 instead).
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Populate the new row 2 with the latest price entry
$ws.Cells.Item(2, 1).Value = 3
$ws.Cells.Item(2, 2).Value = "ALUMINIUM INGOT"
$ws.Cells.Item(2, 3).Value = "IE07"
$ws.Cells.Item(2, 4).Value = 264.35
$ws.Cells.Item(2, 5).Value = "21-08-2025"

# Rebuild hyperlinks so refs line up with the shifted rows
$ws.Hyperlinks.Delete()

$url2 = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$url3 = "https://nalcoindia.com/wp-content/uploads/2025/08/INGOT-15-08-2025.pdf"
$url4 = "https://nalcoindia.com/wp-content/uploads/2019/01/Ingot-07-08-2025.pdf"

$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), $url2, "", "", $url2)
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), $url3, "", "", $url3)
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), $url4, "", "", $url4)

# Hyperlinks.Add stamps the built-in "Hyperlink" style (underline/blue) on
# the cell; the source sheet keeps the plain data-row look instead, so
# restore it from the neighbouring column which was untouched.
$ws.Range("E2:E4").Copy()
$ws.Range("F2:F4").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

$wb.Save()
